$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.108731150627136
$ws.Range("B1").Value = 2.410696506500244
$ws.Range("C1").Value = 5.056607246398926
$ws.Range("D1").Value = 2.28563928604126
$ws.Range("E1").Value = 1.275812268257141
